$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gruppen")

# Change the "Activ" (column A) value from "Yes" to "No" for the rows
# that list groups which should no longer be created by default.
$rows = @(6, 8, 9, 10, 11, 12, 13, 14, 15)
foreach ($r in $rows) {
    $ws.Range("A$r").Value = "No"
}

# Record the active cell/selection as it was when the workbook was saved.
$ws.Range("A24").Select()
